$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its literal text formatting (e.g. "1.00", "11.10", "68.361.23")
# instead of being auto-coerced into a number by Excel when the new values are assigned below.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D19", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D44", "D45", "D46", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.361.23'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '3.364.94'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '585.15'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '178.22'
$ws.Range("E6").Value = '  +1.38%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.186'
$ws.Range("E9").Value = '  +4.25%  '
$ws.Range("D10").Value = '0.583'
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("D11").Value = '48.16'
$ws.Range("E11").Value = '  +6.12%  '
$ws.Range("D12").Value = '0.0000275'
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("D14").Value = '3.909.54'
$ws.Range("E14").Value = '  +0.79%  '
$ws.Range("D15").Value = '8.46'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '68.420.14'
$ws.Range("E16").Value = '  +1.55%  '
$ws.Range("D18").Value = '3.361.84'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").Value = '17.52'
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("E20").Value = '  +2.70%  '
$ws.Range("D21").Value = '0.898'
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").Value = '5.47'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").Value = '17.03'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").Value = '100.51'
$ws.Range("E25").Value = '  +1.96%  '
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("D27").Value = '9.53'
$ws.Range("E27").Value = '  +2.98%  '
$ws.Range("D28").Value = '33.16'
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("D29").Value = '8.56'
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("E30").Value = '  -5.23%  '
$ws.Range("D31").Value = '11.10'
$ws.Range("E31").Value = '  +1.23%  '
$ws.Range("D32").Value = '554.83'
$ws.Range("E32").Value = '  -3.68%  '
$ws.Range("D33").Value = '0.106'
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("D34").Value = '58.08'
$ws.Range("E34").Value = '  +2.58%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '3.723.56'
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("D38").Value = '0.137'
$ws.Range("E38").Value = '  +4.50%  '
$ws.Range("D39").Value = '34.88'
$ws.Range("E39").Value = '  +2.24%  '
$ws.Range("E40").Value = '  +2.40%  '
$ws.Range("D41").Value = '2.63'
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  +1.54%  '
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").Value = '3.25'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("D45").Value = '0.0412'
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("D46").Value = '2.65'
$ws.Range("E46").Value = '  +1.88%  '
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").Value = '131.96'
$ws.Range("E50").Value = '  +2.58%  '
$ws.Range("D51").Value = '2.59'
$ws.Range("E51").Value = '  -1.18%  '
